$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il34"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.018404666666666
$ws.Range("H2").Value = 12.055214
$ws.Range("I2").Value = 0.1784894308593523
$ws.Range("J2").Value = 0.1784894308593523
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.096771
$ws.Range("N2").Value = 0.290313
$ws.Range("O2").Value = 0.02354936372111631
$ws.Range("P2").Value = 0.02354936372111631
$ws.Range("Q2").Value = 0.3888650379979999
$ws.Range("R2").Value = 3.499785341982
$ws.Range("S2").Value = 0.00420331252768193
$ws.Range("T2").Value = 0.00420331252768193

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il34"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.018404666666666
$ws.Range("H3").Value = 12.055214
$ws.Range("I3").Value = 0.1784894308593523
$ws.Range("J3").Value = 0.1784894308593523
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.05416133333333333
$ws.Range("N3").Value = 0.162484
$ws.Range("O3").Value = 0.01318023931019921
$ws.Range("P3").Value = 0.01318023931019921
$ws.Range("Q3").Value = 0.2176421546195555
$ws.Range("R3").Value = 1.958779391576
$ws.Range("S3").Value = 0.00235253341306752
$ws.Range("T3").Value = 0.002352533413067519

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il34"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.018404666666666
$ws.Range("H4").Value = 12.055214
$ws.Range("I4").Value = 0.1784894308593523
$ws.Range("J4").Value = 0.1784894308593523
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.958350666666666
$ws.Range("N4").Value = 11.875052
$ws.Range("O4").Value = 0.9632703969686844
$ws.Range("P4").Value = 0.9632703969686844
$ws.Range("Q4").Value = 15.90625479123644
$ws.Range("R4").Value = 143.156293121128
$ws.Range("S4").Value = 0.1719335849186029
$ws.Range("T4").Value = 0.1719335849186029

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il34"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.352037
$ws.Range("H5").Value = 13.056111
$ws.Range("I5").Value = 0.1933087062267439
$ws.Range("J5").Value = 0.1933087062267439
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.096771
$ws.Range("N5").Value = 0.290313
$ws.Range("O5").Value = 0.02354936372111631
$ws.Range("P5").Value = 0.02354936372111631
$ws.Range("Q5").Value = 0.421150972527
$ws.Range("R5").Value = 3.790358752743
$ws.Range("S5").Value = 0.004552297033392013
$ws.Range("T5").Value = 0.004552297033392013

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il34"
$ws.Range("C6").Value = "Ptprz1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.352037
$ws.Range("H6").Value = 13.056111
$ws.Range("I6").Value = 0.1933087062267439
$ws.Range("J6").Value = 0.1933087062267439
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05416133333333333
$ws.Range("N6").Value = 0.162484
$ws.Range("O6").Value = 0.01318023931019921
$ws.Range("P6").Value = 0.01318023931019921
$ws.Range("Q6").Value = 0.235712126636
$ws.Range("R6").Value = 2.121409139724
$ws.Range("S6").Value = 0.00254785500881348
$ws.Range("T6").Value = 0.00254785500881348

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il34"
$ws.Range("C7").Value = "Ptprz1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.352037
$ws.Range("H7").Value = 13.056111
$ws.Range("I7").Value = 0.1933087062267439
$ws.Range("J7").Value = 0.1933087062267439
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.958350666666666
$ws.Range("N7").Value = 11.875052
$ws.Range("O7").Value = 0.9632703969686844
$ws.Range("P7").Value = 0.9632703969686844
$ws.Range("Q7").Value = 17.226888560308
$ws.Range("R7").Value = 155.041997042772
$ws.Range("S7").Value = 0.1862085541845384
$ws.Range("T7").Value = 0.1862085541845384

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Il34"
$ws.Range("C8").Value = "Ptprz1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.14296233333333
$ws.Range("H8").Value = 42.428887
$ws.Range("I8").Value = 0.6282018629139038
$ws.Range("J8").Value = 0.6282018629139038
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.096771
$ws.Range("N8").Value = 0.290313
$ws.Range("O8").Value = 0.02354936372111631
$ws.Range("P8").Value = 0.02354936372111631
$ws.Range("Q8").Value = 1.368628607959
$ws.Range("R8").Value = 12.317657471631
$ws.Range("S8").Value = 0.01479375416004237
$ws.Range("T8").Value = 0.01479375416004237

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Il34"
$ws.Range("C9").Value = "Ptprz1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.14296233333333
$ws.Range("H9").Value = 42.428887
$ws.Range("I9").Value = 0.6282018629139038
$ws.Range("J9").Value = 0.6282018629139038
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.05416133333333333
$ws.Range("N9").Value = 0.162484
$ws.Range("O9").Value = 0.01318023931019921
$ws.Range("P9").Value = 0.01318023931019921
$ws.Range("Q9").Value = 0.7660016972564443
$ws.Range("R9").Value = 6.894015275307999
$ws.Range("S9").Value = 0.00827985088831821
$ws.Range("T9").Value = 0.008279850888318209

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Il34"
$ws.Range("C10").Value = "Ptprz1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 14.14296233333333
$ws.Range("H10").Value = 42.428887
$ws.Range("I10").Value = 0.6282018629139038
$ws.Range("J10").Value = 0.6282018629139038
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.958350666666666
$ws.Range("N10").Value = 11.875052
$ws.Range("O10").Value = 0.9632703969686844
$ws.Range("P10").Value = 0.9632703969686844
$ws.Range("Q10").Value = 55.98280438079154
$ws.Range("R10").Value = 503.8452394271239
$ws.Range("S10").Value = 0.6051282578655431
$ws.Range("T10").Value = 0.6051282578655431
